$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format Price/Volume columns as Text so numeric-looking strings
# (e.g. "1.000", "98.90") are preserved exactly instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = '25.553.10'
$ws.Range("E2").Value = '  +2.67%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.668.21'
$ws.Range("E3").Value = '  +2.00%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").Value = '  +0.25%  '

# Row 5: BNB
$ws.Range("D5").Value = '236.67'
$ws.Range("E5").Value = '  +0.56%  '

# Row 6: USDC
$ws.Range("E6").Value = '  +0.06%  '

# Row 7: XRP
$ws.Range("D7").Value = '0.4766'
$ws.Range("E7").Value = '  +0.58%  '

# Row 8: Cardano
$ws.Range("D8").Value = '0.2608'
$ws.Range("E8").Value = '  +1.86%  '

# Row 9: Dogecoin
$ws.Range("D9").Value = '0.06167'
$ws.Range("E9").Value = '  +1.61%  '

# Row 10: WrappedEther
$ws.Range("D10").Value = '1.668.31'
$ws.Range("E10").Value = '  +2.05%  '

# Row 11: TRON
$ws.Range("D11").Value = '0.07008'
$ws.Range("E11").Value = '  +0.87%  '

# Row 12: Solana
$ws.Range("D12").Value = '14.78'
$ws.Range("E12").Value = '  +0.59%  '

# Row 13: Polygon
$ws.Range("D13").Value = '0.5872'
$ws.Range("E13").Value = '  -3.67%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '4.365'
$ws.Range("E14").Value = '  +0.83%  '

# Row 15: Litecoin
$ws.Range("D15").Value = '75.28'
$ws.Range("E15").Value = '  +3.27%  '

# Row 16: Dai
$ws.Range("E16").Value = '  +0.02%  '

# Row 17: BinanceUSD
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  +0.25%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '25.553.06'
$ws.Range("E18").Value = '  +2.66%  '

# Row 19: ShibaInu
$ws.Range("D19").Value = '0.000006731'
$ws.Range("E19").Value = '  +2.30%  '

# Row 20: Avalanche
$ws.Range("E20").Value = '  +2.78%  '

# Row 21: WrappedliquidstakedEther2.0
$ws.Range("D21").Value = '1.884.10'
$ws.Range("E21").Value = '  +2.16%  '

# Row 22: Uniswap
$ws.Range("D22").Value = '4.436'
$ws.Range("E22").Value = '  +1.89%  '

# Row 23: Cosmos
$ws.Range("D23").Value = '8.785'
$ws.Range("E23").Value = '  +2.61%  '

# Row 24: Chainlink
$ws.Range("D24").Value = '5.249'
$ws.Range("E24").Value = '  +0.29%  '

# Row 25: Monero
$ws.Range("D25").Value = '136.95'
$ws.Range("E25").Value = '  +2.62%  '

# Row 26: EthereumClassic
$ws.Range("D26").Value = '15.01'
$ws.Range("E26").Value = '  +1.72%  '

# Row 27: Toncoin
$ws.Range("D27").Value = '1.381'
$ws.Range("E27").Value = '  +0.80%  '

# Row 28: LidoDAOToken
$ws.Range("D28").Value = '1.716'
$ws.Range("E28").Value = '  +5.02%  '

# Row 29: BitcoinCash
$ws.Range("D29").Value = '104.69'
$ws.Range("E29").Value = '  +1.83%  '

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").Value = '3.995'
$ws.Range("E30").Value = '  +6.45%  '

# Row 31: Stellar
$ws.Range("D31").Value = '0.07870'
$ws.Range("E31").Value = '  +1.78%  '

# Row 32: Filecoin
$ws.Range("D32").Value = '3.623'
$ws.Range("E32").Value = '  +2.26%  '

# Row 33: Frax
$ws.Range("B33").Value = 'Frax'
$ws.Range("C33").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D33").Value = '0.9993'
$ws.Range("E33").Value = '  +0.07%  '

# Row 34: Hedera
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.04311'
$ws.Range("E34").Value = '  +0.45%  '

# Row 35: HuobiToken
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.624'
$ws.Range("E35").Value = '  +1.14%  '

# Row 36: ARBITRUM
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '0.9536'
$ws.Range("E36").Value = '  +3.59%  '

# Row 37: ImmutableX
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.6046'
$ws.Range("E37").Value = '  +4.40%  '

# Row 38: TrustWalletToken
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '0.9467'
$ws.Range("E38").Value = '  +15.90%  '

# Row 39: MXToken
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.562'
$ws.Range("E39").Value = '  +0.74%  '

# Row 40: PaxDollar
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '0.9998'
$ws.Range("E40").Value = '  +0.19%  '

# Row 41: RenderToken
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '1.847'
$ws.Range("E41").Value = '  +4.04%  '

# Row 42: VeChain
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.01472'
$ws.Range("E42").Value = '  -4.34%  '

# Row 43: Quant
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '98.90'
$ws.Range("E43").Value = '  +1.72%  '

# Row 44: TheSandbox
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.3749'
$ws.Range("E44").Value = '  +1.61%  '

# Row 45: FraxShare
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '4.874'
$ws.Range("E45").Value = '  +3.55%  '

# Row 46: Algorand
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1116'
$ws.Range("E46").Value = '  +2.30%  '

# Row 47: Aptos
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '6.207'
$ws.Range("E47").Value = '  +2.87%  '

# Row 48: Cronos
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05265'
$ws.Range("E48").Value = '  +1.37%  '

# Row 49: Elrond
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '29.91'
$ws.Range("E49").Value = '  +1.61%  '

# Row 50: EnergySwap
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '7.473'
$ws.Range("E50").Value = '  +4.32%  '

# Row 51: TrueUSD
$ws.Range("B51").Value = 'TrueUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D51").Value = '1.002'
$ws.Range("E51").Value = '  +0.27%  '
